# Natmi following Dr Hou advice
# Update LR-pair summary (Bmp6 -> Acvr1) to include the "ECs" sending cluster
# alongside "FAPs" and "sCs", expanding the 2 x 3 cluster-pair grid to 3 x 3
# (rows 2-10) with recomputed specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp6"
$ws.Range("C2").Value = "Acvr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 4.145393333333334
$ws.Range("H2").Value = 12.43618
$ws.Range("I2").Value = 0.1621900462138432
$ws.Range("J2").Value = 0.1621900462138432
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.695610666666666
$ws.Range("N2").Value = 14.086832
$ws.Range("O2").Value = 0.1802066564018305
$ws.Range("P2").Value = 0.1802066564018305
$ws.Range("Q2").Value = 19.46515315352889
$ws.Range("R2").Value = 175.18637838176
$ws.Range("S2").Value = 0.02922772592985506
$ws.Range("T2").Value = 0.02922772592985505

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp6"
$ws.Range("C3").Value = "Acvr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 4.145393333333334
$ws.Range("H3").Value = 12.43618
$ws.Range("I3").Value = 0.1621900462138432
$ws.Range("J3").Value = 0.1621900462138432
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 15.51448033333333
$ws.Range("N3").Value = 46.543441
$ws.Range("O3").Value = 0.5954098039960916
$ws.Range("P3").Value = 0.5954098039960916
$ws.Range("Q3").Value = 64.31362334393111
$ws.Range("R3").Value = 578.82261009538
$ws.Range("S3").Value = 0.09656954362630142
$ws.Range("T3").Value = 0.09656954362630141

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Bmp6"
$ws.Range("C4").Value = "Acvr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 4.145393333333334
$ws.Range("H4").Value = 12.43618
$ws.Range("I4").Value = 0.1621900462138432
$ws.Range("J4").Value = 0.1621900462138432
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.846719333333333
$ws.Range("N4").Value = 17.540158
$ws.Range("O4").Value = 0.2243835396020779
$ws.Range("P4").Value = 0.2243835396020779
$ws.Range("Q4").Value = 24.23695134627111
$ws.Range("R4").Value = 218.13256211644
$ws.Range("S4").Value = 0.03639277665768674
$ws.Range("T4").Value = 0.03639277665768673

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp6"
$ws.Range("C5").Value = "Acvr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.772999666666666
$ws.Range("H5").Value = 20.318999
$ws.Range("I5").Value = 0.2649961151116367
$ws.Range("J5").Value = 0.2649961151116367
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 4.695610666666666
$ws.Range("N5").Value = 14.086832
$ws.Range("O5").Value = 0.1802066564018305
$ws.Range("P5").Value = 0.1802066564018305
$ws.Range("Q5").Value = 31.80336948012977
$ws.Range("R5").Value = 286.230325321168
$ws.Range("S5").Value = 0.04775406386374263
$ws.Range("T5").Value = 0.04775406386374263

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Bmp6"
$ws.Range("C6").Value = "Acvr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6.772999666666666
$ws.Range("H6").Value = 20.318999
$ws.Range("I6").Value = 0.2649961151116367
$ws.Range("J6").Value = 0.2649961151116367
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.51448033333333
$ws.Range("N6").Value = 46.543441
$ws.Range("O6").Value = 0.5954098039960916
$ws.Range("P6").Value = 0.5954098039960916
$ws.Range("Q6").Value = 105.0795701261732
$ws.Range("R6").Value = 945.716131135559
$ws.Range("S6").Value = 0.1577812849583453
$ws.Range("T6").Value = 0.1577812849583453

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Bmp6"
$ws.Range("C7").Value = "Acvr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6.772999666666666
$ws.Range("H7").Value = 20.318999
$ws.Range("I7").Value = 0.2649961151116367
$ws.Range("J7").Value = 0.2649961151116367
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 5.846719333333333
$ws.Range("N7").Value = 17.540158
$ws.Range("O7").Value = 0.2243835396020779
$ws.Range("P7").Value = 0.2243835396020779
$ws.Range("Q7").Value = 39.59982809576022
$ws.Range("R7").Value = 356.3984528618419
$ws.Range("S7").Value = 0.05946076628954872
$ws.Range("T7").Value = 0.05946076628954872

# Row 8: sCs -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Bmp6"
$ws.Range("C8").Value = "Acvr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.640471
$ws.Range("H8").Value = 43.921413
$ws.Range("I8").Value = 0.5728138386745201
$ws.Range("J8").Value = 0.5728138386745202
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 4.695610666666666
$ws.Range("N8").Value = 14.086832
$ws.Range("O8").Value = 0.1802066564018305
$ws.Range("P8").Value = 0.1802066564018305
$ws.Range("Q8").Value = 68.745951792624
$ws.Range("R8").Value = 618.713566133616
$ws.Range("S8").Value = 0.1032248666082328
$ws.Range("T8").Value = 0.1032248666082328

# Row 9: sCs -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Bmp6"
$ws.Range("C9").Value = "Acvr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.640471
$ws.Range("H9").Value = 43.921413
$ws.Range("I9").Value = 0.5728138386745201
$ws.Range("J9").Value = 0.5728138386745202
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 15.51448033333333
$ws.Range("N9").Value = 46.543441
$ws.Range("O9").Value = 0.5954098039960916
$ws.Range("P9").Value = 0.5954098039960916
$ws.Range("Q9").Value = 227.139299400237
$ws.Range("R9").Value = 2044.253694602133
$ws.Range("S9").Value = 0.3410589754114448
$ws.Range("T9").Value = 0.3410589754114449

# Row 10: sCs -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Bmp6"
$ws.Range("C10").Value = "Acvr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.640471
$ws.Range("H10").Value = 43.921413
$ws.Range("I10").Value = 0.5728138386745201
$ws.Range("J10").Value = 0.5728138386745202
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.846719333333333
$ws.Range("N10").Value = 17.540158
$ws.Range("O10").Value = 0.2243835396020779
$ws.Range("P10").Value = 0.2243835396020779
$ws.Range("Q10").Value = 85.59872484480599
$ws.Range("R10").Value = 770.3885236032539
$ws.Range("S10").Value = 0.1285299966548424
$ws.Range("T10").Value = 0.1285299966548424

